# Applies the OOXML change captured by the source diff:
#   ppt/slides/slide16.xml : the <a:tbl> inside the graphicFrame
#   "Google Shape;213;p29" switches its <a:tableStyleId> from the
#   custom "Table_0" style {CFD55213-2E06-4827-8068-D3A21DF76405}
#   (defined in ppt/tableStyles.xml) to the built-in table style
#   {6CD20A6D-EBC6-4B87-AD29-A77E924511F2}.
#
# PowerPoint's Table object does not allow `.Style` to be assigned
# directly (it is read-only through the property setter) — the COM
# surface requires `Table.ApplyStyle("{GUID}")` instead, mirroring
# what the "Table Styles" gallery does in the UI.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

$targetStyleId = "{6CD20A6D-EBC6-4B87-AD29-A77E924511F2}"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle($targetStyleId)
    }
}
